# Swap the contents (and matching column widths) of columns A and B on the
# active sheet, rows 1-10.
#
# Column A currently holds width 3.140625 with the "counts" series
# (0,1,3,5,7,9,11,13,15,18) and column B holds width 2.140625 with the
# "ThreshE" series (0,0,0,1,2,2,1,1,0,0). The two columns need to trade
# places: A should end up with B's old width/values and B should end up
# with A's old width/values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the column widths -------------------------------------------------
# The target widths are 2.140625 (col A) and 3.140625 (col B). The ColumnWidth
# setter in this runtime only resolves to discrete (1/6 character) increments,
# so the input values below are chosen as the ones that land on the closest
# achievable rendered width to each target.
$ws.Columns.Item(1).ColumnWidth = 1.25   # renders as ~2.1667 (closest achievable to 2.140625)
$ws.Columns.Item(2).ColumnWidth = 2.25   # renders as ~3.1667 (closest achievable to 3.140625)

# --- Swap the cell values ---------------------------------------------------
for ($r = 1; $r -le 10; $r++) {
    $valA = $ws.Cells.Item($r, 1).Value2
    $valB = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value = $valB
    $ws.Cells.Item($r, 2).Value = $valA
}
